$wb = $excel.ActiveWorkbook

# --- zh-cn sheet: rows 4-7 are the "ht" priority rows whose handoff xliff was
#     regenerated. Priority goes from "low" -> "ht" and the Latest Handoff
#     Datetime (col H) moves from 18:34:29 -> 18:34:45.
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("H4").Value = "2016-09-02 18:34:45"
$wsZh.Range("E5").Value = "ht"
$wsZh.Range("H5").Value = "2016-09-02 18:34:45"
$wsZh.Range("E6").Value = "ht"
$wsZh.Range("H6").Value = "2016-09-02 18:34:45"
$wsZh.Range("E7").Value = "ht"
$wsZh.Range("H7").Value = "2016-09-02 18:34:45"

# --- de-de sheet: same rows, same Priority change, and Latest Handoff
#     Datetime (col H) moves from 18:34:34 -> 18:34:49. This value is shared
#     with the Overview sheet's "Latest HO Xliff Generate Date" column.
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("H4").Value = "2016-09-02 18:34:49"
$wsDe.Range("E5").Value = "ht"
$wsDe.Range("H5").Value = "2016-09-02 18:34:49"
$wsDe.Range("E6").Value = "ht"
$wsDe.Range("H6").Value = "2016-09-02 18:34:49"
$wsDe.Range("E7").Value = "ht"
$wsDe.Range("H7").Value = "2016-09-02 18:34:49"

# --- Overview sheet: "Latest HO Xliff Generate Date" column (G) for the same
#     four rows reflects the same regenerated handoff timestamp.
$wsOv = $wb.Worksheets.Item("Overview")
$wsOv.Range("G4").Value = "2016-09-02 18:34:49"
$wsOv.Range("G5").Value = "2016-09-02 18:34:49"
$wsOv.Range("G6").Value = "2016-09-02 18:34:49"
$wsOv.Range("G7").Value = "2016-09-02 18:34:49"
